# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - row number => new F value
$exhibitionUpdates = @{
    2  = 3041
    3  = 466
    5  = 33
    6  = 261
    7  = 214
    8  = 14538
    9  = 161
    10 = 123
    11 = 5806
    13 = 73
    15 = 63
    17 = 14
    19 = 183
    20 = 796
    21 = 2940
    22 = 55
    23 = 10592
    24 = 1201
    25 = 61
    26 = 84
    28 = 245
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (all types) - row number => new F value
$allTypesUpdates = @{
    3  = 3041
    4  = 466
    6  = 33
    7  = 261
    8  = 214
    9  = 14538
    10 = 161
    11 = 123
    12 = 5806
    14 = 73
    16 = 63
    18 = 14
    20 = 183
    21 = 796
    22 = 2940
    23 = 55
    25 = 10592
    26 = 1201
    27 = 61
    28 = 84
    30 = 245
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
